$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Remove the stray "1" (ParentID) entries that were mistakenly placed on the
# top-level category rows (B12, B15, B18, B20) ...
$ws.Range("B12").Clear()
$ws.Range("B15").Clear()
$ws.Range("B18").Clear()
$ws.Range("B20").Clear()

# ... and move that value to where it actually belongs: row 27 (B1. Stammdaten
# & Betreiber), which was missing its ParentID.
$ws.Range("B27").Value2 = "1"

# Reflect the author's final on-screen selection/scroll position.
[void]$ws.Activate()
[void]$ws.Range("B28").Select()
